$d = $word.ActiveDocument

# --- helpers -----------------------------------------------------------
#
# `InsertParagraphAfter()` always manifests the brand-new paragraph it
# creates with a placeholder empty run (`<w:r/>`) in the underlying OOXML,
# even though the "donor" paragraph it was called on is left untouched. We
# lean on that asymmetry to get clean results either way:
#
#  * Clean-Paragraph: rewrites an existing (possibly multi-run, spell-check
#    -annotated) paragraph's text as a single plain run. It inserts a new
#    paragraph right after the target, puts the replacement text into that
#    new paragraph (which folds the placeholder run into a normal text run
#    - no stray empty run survives once it holds text), and then deletes
#    the whole of the original paragraph (its text *and* its paragraph
#    mark) - carrying away any w:proofErr spell-check markers that were
#    anchored to it.
#
#  * New-BlankParagraph: inserts a brand-new, fully-empty paragraph
#    (canonical bare `<w:p/>`, not `<w:p><w:r/></w:p>`) right after the
#    given paragraph index. It temporarily types a placeholder word into
#    the freshly-created paragraph (so the placeholder run becomes a real
#    text run) and then clears the text back out again - clearing a run
#    that actually held text removes the run entirely, instead of leaving
#    an empty run behind.

function Clean-Paragraph($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($index + 1)
    $newPara.Range.InsertBefore($newText)
    $old = $d.Paragraphs.Item($index)
    $old.Range.Delete()
}

function New-BlankParagraph($index) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($index + 1)
    $newPara.Range.InsertBefore("x")
    $r = $newPara.Range
    $r.End = $r.End - 1
    $r.Text = ""
}

# --- rewrite the three existing notes as single clean runs -------------
# (processed bottom-to-top so earlier paragraph indices stay valid while
# later ones are being rewritten)

Clean-Paragraph 5 "Boksz model – content, margin, padding, width, height, border"
Clean-Paragraph 3 "Inheritance – nestelt dolgok örökölnek tulajdonságokat"
Clean-Paragraph 1 "Descendant selector – leszármazó szelektor, pl a p-ben lévő span-eket akarom formázni, akkor a kettőt szóközzel választom el."

# --- add the new "Margin" note, preceded by a blank separator line -----
# (paragraph 5 is "Boksz model ..." again, same index as before since
# Clean-Paragraph above is index-count-neutral)

New-BlankParagraph 5
$marginPara = $d.Paragraphs.Item(7)
$marginPara.Range.InsertBefore("Margin: fenn, jobb, lenn, bal (shorthand)")
